# add error message pop-ups to 2_add_layouts/app.R
# Updates the workflow diagram sheet: adds a new "merge_layouts app.R" row
# and a new "Addressed" column documenting error handling, and re-shapes
# the Input/Script/Output table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells whose old content is being removed/relocated so the
# table doesn't retain stray leftovers.
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""

# --- Write new/changed cells, introducing brand-new strings in the
#     order they first appear so the shared-string table lines up. ---
$ws.Range("A3").Value = "formatted data"
$ws.Range("E2").Value = "formatted data (nested by each input data column)"
$ws.Range("B3").Value = "merge_layouts app.R"
$ws.Range("D1").Value = "Addressed"
$ws.Range("D2").Value = "Differences in raw data upload formats"
$ws.Range("E4").Value = "labeled data"

$ws.Range("E1").Value = "Output"
$ws.Range("E3").Value = "Tmas from dRFU"
$ws.Range("A7").Value = "nested formatted data"
$ws.Range("E7").Value = "Tmas from model fitting"
$ws.Range("A8").Value = "layout file"
$ws.Range("E8").Value = "formatted layout "

# --- Column widths (targets: 20.83203125 / 21.6640625 / 34 / 44.5 chars;
#     the values below are tuned so the COM layer's pixel-quantized
#     ColumnWidth setter lands on the closest achievable stored width). ---
$ws.Columns.Item(2).ColumnWidth = 20.001
$ws.Columns.Item(3).ColumnWidth = 20.834
$ws.Columns.Item(4).ColumnWidth = 33.1675
$ws.Columns.Item(5).ColumnWidth = 43.6675

# --- Selection ---
$ws.Range("B13").Select()
